# Tutorial 6 solution update:
#  - Reformat the Date column (A3:A21) from "dd/mm/yyyy" to "dd-mm-yyyy".
#  - Correct the attendance tallies (Total/Real/Invalid/Absent) for the
#    rows where the duplicate/invalid/absent classification changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to Text first so strings like "01-08-2022" are
# kept as literal text instead of being auto-parsed into a date serial.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 1).Value = $dates[$row]
}

# Updated attendance counts: Total(D), Real(E), Invalid(G), Absent(H)
# Only rows 3, 5, 10 and 12 change numerically.
$counts = @{
    3  = @{ D = 1; E = 0; G = 1; H = 1 }
    5  = @{ D = 1; E = 1; G = 0; H = 0 }
    10 = @{ D = 1; E = 1; G = 0; H = 0 }
    12 = @{ D = 1; E = 1; G = 0; H = 0 }
}

foreach ($row in $counts.Keys) {
    $rowVals = $counts[$row]
    $ws.Cells.Item($row, 4).Value = $rowVals.D
    $ws.Cells.Item($row, 5).Value = $rowVals.E
    $ws.Cells.Item($row, 7).Value = $rowVals.G
    $ws.Cells.Item($row, 8).Value = $rowVals.H
}
